# Update "想去人数" (number of people interested) figures that changed
# between the two generated-data snapshots.
#
# Sheet "展览" (sheet1): F3 1956->1968, F4 841->842, F5 984->1006, F6 344->345
# Sheet "全部类型" (sheet4): F3 1956->1968, F5 841->842, F6 984->1006, F7 345

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1968
$ws1.Range("F4").Value = 842
$ws1.Range("F5").Value = 1006
$ws1.Range("F6").Value = 345

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1968
$ws4.Range("F5").Value = 842
$ws4.Range("F6").Value = 1006
$ws4.Range("F7").Value = 345
